# Change format of dates
#
# Summary of the edit (per commit "Change format of dates"):
#  - The "task" sheet's dueDate column (G) becomes "dueDateTime": it now stores
#    the due moment as a date+time value (end of day, 23:59) instead of a bare
#    date, and is displayed with a new custom format "m/d/yy h:mm;@".
#  - The "date" column (K) keeps its values but switches to a new custom date
#    format "m/d/yy;@" (instead of the built-in "m/d/yyyy").
#  - A block of additional rows below the existing data was pre-formatted with
#    the new date/date-time formats (ready for future rows to be filled in).
#  - The used range / selection grew accordingly.

$wb = $excel.ActiveWorkbook
$taskWs = $wb.Worksheets.Item(1)
$userWs = $wb.Worksheets.Item(2)

# --- "task" sheet -----------------------------------------------------

# Rename the header of column G from "dueDate" to "dueDateTime"
$taskWs.Cells.Item(1, 7).Value = "dueDateTime"

# Convert the existing dueDate values (date-only) to end-of-day timestamps
# and apply the new date-time display format.
$taskWs.Cells.Item(2, 7).Value = 43070.999305555553
$taskWs.Cells.Item(2, 7).NumberFormat = "m/d/yy\ h:mm;@"

$taskWs.Cells.Item(3, 7).Value = 43074.999305555553
$taskWs.Cells.Item(3, 7).NumberFormat = "m/d/yy\ h:mm;@"

# The "date" column (K) keeps its values, only the display format changes.
$taskWs.Cells.Item(2, 11).NumberFormat = "m/d/yy;@"
$taskWs.Cells.Item(3, 11).NumberFormat = "m/d/yy;@"

# Pre-format additional (currently empty) rows below the data with the new
# date/date-time formats, mirroring the formatting carried down the column.
$taskWs.Range($taskWs.Cells.Item(4, 7), $taskWs.Cells.Item(9, 7)).NumberFormat = "m/d/yy\ h:mm;@"
$taskWs.Range($taskWs.Cells.Item(4, 11), $taskWs.Cells.Item(16, 11)).NumberFormat = "m/d/yy;@"

# Update the dimension / selection to reflect the newly used range.
$taskWs.Range("O21").Select()

# --- "user" sheet -------------------------------------------------------
# No data changed on this sheet; the underlying shared-string reshuffle
# caused by the dueDate -> dueDateTime rename is handled automatically by
# the engine when the strings are written/read through the object model.
$userWs.Range("B6").Select()
$taskWs.Activate()
